# "run prepare & render with final data"
# Adds a Russia data column (inserted between Japan and Saudi Arabia),
# refreshes row labels for two survey questions, and writes the final
# recomputed percentages for every country/question cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column for Russia, right after Japan (column K) ---
# This shifts the old Saudi Arabia column (L) to M and the old USA
# column (M) to N, preserving their existing values/blanks.
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L1").Value = "Russia"

# --- Row-label text refresh ---
$ws.Range("A2").Value = "Supports tax on world top 1% to finance global poverty reduction`n(Additional 15% tax on income over [`$120k/year in PPP])"
$ws.Range("A3").Value = "Supports tax on world top 3% to finance global poverty reduction`n(Additional 15% tax over [`$80k], 30% over [`$120k], 45% over [`$1M])"
$ws.Range("A4").Value = "Prefers sustainable future"
$ws.Range("A5").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""
$ws.Range("A6").Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"
$ws.Range("A7").Value = "More likely to vote for party if part of worldwide`ncoalition for climate action and global redistribution"
$ws.Range("A8").Value = "Supports reparations for colonization and slavery in`nthe form of funding education and technology transfers"
$ws.Range("A9").Value = "`"My taxes should go towards solving global problems`""

# --- Final recomputed data (columns B:N, rows 2:9) ---
# Row 2
$ws.Range("B2").Value = 0.537951467268623
$ws.Range("C2").Value = 0.588318584070797
$ws.Range("D2").Value = 0.597746584709126
$ws.Range("E2").Value = 0.614782732381824
$ws.Range("F2").Value = 0.699507350357918
$ws.Range("G2").Value = 0.492647690680722
$ws.Range("H2").Value = 0.58648437760578
$ws.Range("I2").Value = 0.547405768931041
$ws.Range("J2").Value = 0.524553235889069
$ws.Range("K2").Value = 0.389856557970562
$ws.Range("L2").Value = 0.582306191819793
$ws.Range("M2").Value = 0.68442794362766
$ws.Range("N2").Value = 0.490668612212663
# Row 3
$ws.Range("B3").Value = 0.492862092862093
$ws.Range("C3").Value = 0.535241502683363
$ws.Range("D3").Value = 0.578029654074004
$ws.Range("E3").Value = 0.531972070523522
$ws.Range("F3").Value = 0.596746035930078
$ws.Range("G3").Value = 0.522317740589736
$ws.Range("H3").Value = 0.552425813190738
$ws.Range("I3").Value = 0.556404416090476
$ws.Range("J3").Value = 0.336017347574206
$ws.Range("K3").Value = 0.32128421797297
$ws.Range("L3").Value = 0.591124471229443
$ws.Range("M3").Value = 0.666290156760115
$ws.Range("N3").Value = 0.448211203307958
# Row 4
$ws.Range("B4").Value = 0.664778686805119
$ws.Range("C4").Value = 0.680960854092527
$ws.Range("D4").Value = 0.703421739255081
$ws.Range("E4").Value = 0.68701219464814
$ws.Range("F4").Value = 0.72909694938135
$ws.Range("G4").Value = 0.569764105372422
$ws.Range("H4").Value = 0.726504855911
$ws.Range("I4").Value = 0.672412949788013
$ws.Range("J4").Value = 0.653357988317707
$ws.Range("K4").Value = 0.706011505764104
$ws.Range("L4").Value = 0.6876153744808
$ws.Range("M4").Value = 0.666711263875173
$ws.Range("N4").Value = 0.603123607121577
# Row 5
$ws.Range("B5").Value = 0.553737500874065
$ws.Range("C5").Value = 0.611565836298932
$ws.Range("D5").Value = 0.576391524303972
$ws.Range("E5").Value = 0.6111956659185
$ws.Range("F5").Value = 0.683237385490521
$ws.Range("G5").Value = 0.681995327066417
$ws.Range("H5").Value = 0.698497778774273
$ws.Range("I5").Value = 0.51480902544306
$ws.Range("J5").Value = 0.535091205491473
$ws.Range("K5").Value = 0.498450958769852
$ws.Range("L5").Value = 0.766211423943959
$ws.Range("M5").Value = 0.572921342842471
$ws.Range("N5").Value = 0.426115190401503
# Row 6
$ws.Range("B6").Value = 0.611777124330845
$ws.Range("C6").Value = 0.672953736654804
$ws.Range("D6").Value = 0.668600432112831
$ws.Range("E6").Value = 0.673402016569035
$ws.Range("F6").Value = 0.737070007431593
$ws.Range("G6").Value = 0.655329158659857
$ws.Range("H6").Value = 0.724297246090431
$ws.Range("I6").Value = 0.645148556496296
$ws.Range("J6").Value = 0.602999620241999
$ws.Range("K6").Value = 0.507428712494617
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value = 0.55354536542297
$ws.Range("N6").Value = 0.612788981791228
# Row 7
$ws.Range("B7").Value = 0.350712198685172
$ws.Range("C7").Value = 0.39870648658931
$ws.Range("D7").Value = 0.431546204389394
$ws.Range("E7").Value = 0.386040021908065
$ws.Range("F7").Value = 0.479597303697709
$ws.Range("G7").Value = 0.287167476601508
$ws.Range("H7").Value = 0.43772759060628
$ws.Range("I7").Value = 0.395885893816819
$ws.Range("J7").Value = 0.325907437125978
$ws.Range("K7").Value = 0.217091601807433
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0.370801286983423
# Row 8
$ws.Range("B8").Value = 0.351261177098625
$ws.Range("C8").Value = 0.384560570071259
$ws.Range("D8").Value = 0.326868878117742
$ws.Range("E8").Value = 0.334505390580585
$ws.Range("F8").Value = 0.535274410980465
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0.401939935801779
$ws.Range("I8").Value = 0.348746332347973
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 0.300404564269368
# Row 9
$ws.Range("B9").Value = 0.424604085379849
$ws.Range("C9").Value = 0.437511888910025
$ws.Range("D9").Value = 0.328866922959613
$ws.Range("E9").Value = 0.453321796373249
$ws.Range("F9").Value = 0.513351118669645
$ws.Range("G9").Value = 0.41080678800161
$ws.Range("H9").Value = 0.500387925343261
$ws.Range("I9").Value = 0.450348125937886
$ws.Range("J9").Value = 0.388060989658246
$ws.Range("K9").Value = 0.314322752209368
$ws.Range("L9").Value = 0.400915137390113
$ws.Range("M9").Value = 0.653477703137907
$ws.Range("N9").Value = 0.406101165624595
